$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9996811035619473
$ws.Range("C2").Value = 1.056680008865375
$ws.Range("D2").Value = 1.535658889322983
$ws.Range("E2").Value = 1.239217046898155
$ws.Range("F2").Value = 0.7599672122005358

$ws.Range("B3").Value = 0.6541542769893847
$ws.Range("C3").Value = 0.8341207046287089
$ws.Range("D3").Value = 0.8907909344799998
$ws.Range("E3").Value = 0.9438172145495122
$ws.Range("F3").Value = 0.7171495399430653
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = 0.7114172426094275
$ws.Range("C4").Value = 0.7114172426094275
$ws.Range("D4").Value = 0.7732681336315609
$ws.Range("E4").Value = 0.8793566589453684
$ws.Range("F4").Value = 0.5662017031583992
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 0.9993615165702219
$ws.Range("C5").Value = 0.9993615165702219
$ws.Range("D5").Value = 1.585080129446412
$ws.Range("E5").Value = 1.258999654267789
$ws.Range("F5").Value = 1.082918915380905
$ws.Range("G5").Value = 2

$ws.Rows("6:9").Delete()
